# Update cryptocurrency price (D) and 1h volume-change (E) columns
# with the latest scraped values from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ D = "71.704.44"; E = "  +3.19%  " }
    3 = @{ D = "3.586.52"; E = "  +5.71%  " }
    4 = @{ D = "0.996"; E = "  -0.33%  " }
    5 = @{ D = "593.05"; E = "  +0.92%  " }
    6 = @{ D = "183.17"; E = "  +1.60%  " }
    7 = @{ D = "3.568.10"; E = "  +5.32%  " }
    8 = @{ D = "0.607"; E = "  +1.86%  " }
    9 = @{ D = "0.999"; E = "  -0.05%  " }
    10 = @{ D = "0.207"; E = "  +6.11%  " }
    11 = @{ D = "0.607"; E = "  +2.40%  " }
    12 = @{ D = "50.15"; E = "  +3.02%  " }
    13 = @{ D = "0.0000289"; E = "  +2.36%  " }
    14 = @{ D = "695.45"; E = "  +2.00%  " }
    15 = @{ D = "4.139.91"; E = "  +5.18%  " }
    16 = @{ D = "8.90"; E = "  +2.89%  " }
    17 = @{ D = "71.408.59"; E = "  +2.72%  " }
    18 = @{ D = "3.519.10"; E = "  +3.68%  " }
    19 = @{ E = "  +1.41%  " }
    20 = @{ D = "18.30"; E = "  +3.39%  " }
    21 = @{ D = "11.69"; E = "  +3.40%  " }
    22 = @{ D = "0.929"; E = "  +2.60%  " }
    23 = @{ D = "5.54"; E = "  +2.24%  " }
    24 = @{ D = "17.62"; E = "  +2.34%  " }
    25 = @{ D = "104.26"; E = "  +0.32%  " }
    26 = @{ D = "4.02"; E = "  +2.12%  " }
    27 = @{ D = "2.80"; E = "  +2.34%  " }
    28 = @{ D = "10.01"; E = "  +2.35%  " }
    29 = @{ E = "  +3.24%  " }
    30 = @{ D = "9.06"; E = "  +3.60%  " }
    31 = @{ D = "7.49"; E = "  +6.15%  " }
    32 = @{ D = "4.14"; E = "  +15.78%  " }
    33 = @{ D = "587.79"; E = "  +5.08%  " }
    34 = @{ D = "11.30"; E = "  +0.89%  " }
    35 = @{ D = "0.107"; E = "  -0.04%  " }
    36 = @{ D = "59.74"; E = "  +2.50%  " }
    37 = @{ E = "  -0.06%  " }
    38 = @{ D = "3.674.85"; E = "  -0.57%  " }
    39 = @{ D = "0.145"; E = "  +4.18%  " }
    40 = @{ D = "36.06"; E = "  +1.25%  " }
    41 = @{ D = "0.0₃0770"; E = "  +9.87%  " }
    42 = @{ D = "3.47"; E = "  +5.76%  " }
    43 = @{ D = "2.82"; E = "  +3.38%  " }
    44 = @{ D = "0.0439"; E = "  +3.46%  " }
    45 = @{ D = "0.347"; E = "  +1.90%  " }
    46 = @{ D = "3.39"; E = "  +2.26%  " }
    47 = @{ D = "2.76"; E = "  +2.69%  " }
    48 = @{ E = "  +5.06%  " }
    50 = @{ E = "  -0.25%  " }
    51 = @{ D = "133.80"; E = "  +0.70%  " }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    if ($u.ContainsKey("D")) {
        $priceCell = $ws.Range("D$row")
        # Force the price to be stored as text (it may look like a number,
        # e.g. "0.996" or "593.05") so it keeps the same inline-string type
        # as the rest of the sheet instead of being coerced to a number.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.D
        $priceCell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u.E
    }
}
